$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  24"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# Crime data table rows 14-30 (column A text labels are unchanged; only
# shared-string indices shift internally, which the engine manages itself)
# Row 14
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "0"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "***.*"
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = "0"
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = -100
$ws.Cells.Item(14, 9).Value = 2
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(14, 14).Value = -77.777777777777

# Row 15
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "***.*"
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 100
$ws.Cells.Item(15, 9).Value = 11
$ws.Cells.Item(15, 10).Value = 8
$ws.Cells.Item(15, 11).Value = 37.5
$ws.Cells.Item(15, 12).Value = 22.222222222222
$ws.Cells.Item(15, 13).Value = 57.142857142857
$ws.Cells.Item(15, 14).Value = -21.428571428571

# Row 16
$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 4).Value = 6
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 23
$ws.Cells.Item(16, 7).Value = 18
$ws.Cells.Item(16, 8).Value = 27.777777777777
$ws.Cells.Item(16, 9).Value = 100
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 25
$ws.Cells.Item(16, 12).Value = 96.078431372549
$ws.Cells.Item(16, 13).Value = -9.090909090909
$ws.Cells.Item(16, 14).Value = -77.52808988764

# Row 17
$ws.Cells.Item(17, 3).Value = 13
$ws.Cells.Item(17, 4).Value = 9
$ws.Cells.Item(17, 5).Value = 44.444444444444
$ws.Cells.Item(17, 6).Value = 27
$ws.Cells.Item(17, 7).Value = 28
$ws.Cells.Item(17, 8).Value = -3.571428571428
$ws.Cells.Item(17, 9).Value = 133
$ws.Cells.Item(17, 10).Value = 129
$ws.Cells.Item(17, 11).Value = 3.100775193798
$ws.Cells.Item(17, 12).Value = 25.471698113207
$ws.Cells.Item(17, 13).Value = 17.699115044247
$ws.Cells.Item(17, 14).Value = 2.307692307692

# Row 18
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 11
$ws.Cells.Item(18, 5).Value = -72.727272727272
$ws.Cells.Item(18, 6).Value = 6
$ws.Cells.Item(18, 7).Value = 25
$ws.Cells.Item(18, 8).Value = -76
$ws.Cells.Item(18, 9).Value = 79
$ws.Cells.Item(18, 10).Value = 120
$ws.Cells.Item(18, 11).Value = -34.166666666666
$ws.Cells.Item(18, 12).Value = -28.181818181818
$ws.Cells.Item(18, 13).Value = -57.52688172043
$ws.Cells.Item(18, 14).Value = -91.459459459459

# Row 19
$ws.Cells.Item(19, 3).Value = 15
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = 25
$ws.Cells.Item(19, 6).Value = 59
$ws.Cells.Item(19, 7).Value = 47
$ws.Cells.Item(19, 8).Value = 25.531914893617
$ws.Cells.Item(19, 9).Value = 302
$ws.Cells.Item(19, 10).Value = 286
$ws.Cells.Item(19, 11).Value = 5.594405594405
$ws.Cells.Item(19, 12).Value = 41.784037558685
$ws.Cells.Item(19, 13).Value = 61.497326203208
$ws.Cells.Item(19, 14).Value = 9.818181818181

# Row 20
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 7
$ws.Cells.Item(20, 5).Value = -42.857142857142
$ws.Cells.Item(20, 6).Value = 29
$ws.Cells.Item(20, 7).Value = 22
$ws.Cells.Item(20, 8).Value = 31.818181818181
$ws.Cells.Item(20, 9).Value = 156
$ws.Cells.Item(20, 10).Value = 115
$ws.Cells.Item(20, 11).Value = 35.652173913043
$ws.Cells.Item(20, 12).Value = 64.210526315789
$ws.Cells.Item(20, 13).Value = -7.142857142857
$ws.Cells.Item(20, 14).Value = -90.791027154663

# Row 21
$ws.Cells.Item(21, 3).Value = 43
$ws.Cells.Item(21, 4).Value = 45
$ws.Cells.Item(21, 5).Value = -4.444444444444
$ws.Cells.Item(21, 6).Value = 146
$ws.Cells.Item(21, 7).Value = 142
$ws.Cells.Item(21, 8).Value = 2.81690140845
$ws.Cells.Item(21, 9).Value = 783
$ws.Cells.Item(21, 10).Value = 740
$ws.Cells.Item(21, 11).Value = 5.81081081081
$ws.Cells.Item(21, 12).Value = 33.617747440273
$ws.Cells.Item(21, 13).Value = 1.424870466321
$ws.Cells.Item(21, 14).Value = -77.577319587628

# Row 22
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "***.*"
$ws.Cells.Item(22, 6).Value = 4
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = 33.333333333333
$ws.Cells.Item(22, 9).Value = 13
$ws.Cells.Item(22, 10).Value = 9
$ws.Cells.Item(22, 11).Value = 44.444444444444
$ws.Cells.Item(22, 12).Value = 550
$ws.Cells.Item(22, 13).Value = 30
$ws.Cells.Item(22, 14).NumberFormat = "@"
$ws.Cells.Item(22, 14).Value = "***.*"

# Row 23
$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = "0"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "***.*"
$ws.Cells.Item(23, 6).NumberFormat = "@"
$ws.Cells.Item(23, 6).Value = "0"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "0"
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = "***.*"
$ws.Cells.Item(23, 9).NumberFormat = "@"
$ws.Cells.Item(23, 9).Value = "0"
$ws.Cells.Item(23, 10).NumberFormat = "@"
$ws.Cells.Item(23, 10).Value = "0"
$ws.Cells.Item(23, 11).NumberFormat = "@"
$ws.Cells.Item(23, 11).Value = "***.*"
$ws.Cells.Item(23, 12).NumberFormat = "@"
$ws.Cells.Item(23, 12).Value = "***.*"
$ws.Cells.Item(23, 13).NumberFormat = "@"
$ws.Cells.Item(23, 13).Value = "***.*"
$ws.Cells.Item(23, 14).NumberFormat = "@"
$ws.Cells.Item(23, 14).Value = "***.*"

# Row 24
$ws.Cells.Item(24, 3).Value = 44
$ws.Cells.Item(24, 4).Value = 19
$ws.Cells.Item(24, 5).Value = 131.578947368421
$ws.Cells.Item(24, 6).Value = 116
$ws.Cells.Item(24, 7).Value = 120
$ws.Cells.Item(24, 8).Value = -3.333333333333
$ws.Cells.Item(24, 9).Value = 580
$ws.Cells.Item(24, 10).Value = 628
$ws.Cells.Item(24, 11).Value = -7.64331210191
$ws.Cells.Item(24, 12).Value = -8.517350157728
$ws.Cells.Item(24, 13).Value = 15.768463073852
$ws.Cells.Item(24, 14).NumberFormat = "@"
$ws.Cells.Item(24, 14).Value = "***.*"

# Row 25
$ws.Cells.Item(25, 3).Value = 8
$ws.Cells.Item(25, 4).Value = 15
$ws.Cells.Item(25, 5).Value = -46.666666666666
$ws.Cells.Item(25, 6).Value = 30
$ws.Cells.Item(25, 7).Value = 51
$ws.Cells.Item(25, 8).Value = -41.176470588235
$ws.Cells.Item(25, 9).Value = 233
$ws.Cells.Item(25, 10).Value = 233
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 18.274111675126
$ws.Cells.Item(25, 13).Value = -32.658959537572
$ws.Cells.Item(25, 14).NumberFormat = "@"
$ws.Cells.Item(25, 14).Value = "***.*"

# Row 26
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "***.*"
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(26, 7).Value = 2
$ws.Cells.Item(26, 8).Value = 50
$ws.Cells.Item(26, 9).Value = 16
$ws.Cells.Item(26, 10).Value = 15
$ws.Cells.Item(26, 11).Value = 6.666666666666
$ws.Cells.Item(26, 12).Value = 45.454545454545
$ws.Cells.Item(26, 13).NumberFormat = "@"
$ws.Cells.Item(26, 13).Value = "***.*"
$ws.Cells.Item(26, 14).NumberFormat = "@"
$ws.Cells.Item(26, 14).Value = "***.*"

# Row 27
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "0"
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 6).Value = 6
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 31
$ws.Cells.Item(27, 10).Value = 31
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 19.230769230769
$ws.Cells.Item(27, 13).NumberFormat = "@"
$ws.Cells.Item(27, 13).Value = "***.*"
$ws.Cells.Item(27, 14).NumberFormat = "@"
$ws.Cells.Item(27, 14).Value = "***.*"

# Row 28
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "0"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "***.*"
$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value = "0"
$ws.Cells.Item(28, 7).Value = 1
$ws.Cells.Item(28, 8).Value = -100
$ws.Cells.Item(28, 9).Value = 2
$ws.Cells.Item(28, 10).Value = 7
$ws.Cells.Item(28, 11).Value = -71.428571428571
$ws.Cells.Item(28, 12).NumberFormat = "@"
$ws.Cells.Item(28, 12).Value = "***.*"
$ws.Cells.Item(28, 13).Value = 100
$ws.Cells.Item(28, 14).Value = -83.333333333333

# Row 29
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "0"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "***.*"
$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value = "0"
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = -100
$ws.Cells.Item(29, 9).Value = 2
$ws.Cells.Item(29, 10).Value = 6
$ws.Cells.Item(29, 11).Value = -66.666666666666
$ws.Cells.Item(29, 12).NumberFormat = "@"
$ws.Cells.Item(29, 12).Value = "***.*"
$ws.Cells.Item(29, 13).Value = 100
$ws.Cells.Item(29, 14).Value = -81.818181818181

# Row 30
$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "0"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "***.*"
$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = "0"
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 8).Value = -100
$ws.Cells.Item(30, 9).NumberFormat = "@"
$ws.Cells.Item(30, 9).Value = "0"
$ws.Cells.Item(30, 10).Value = 4
$ws.Cells.Item(30, 11).Value = -100
$ws.Cells.Item(30, 12).Value = -100
$ws.Cells.Item(30, 13).NumberFormat = "@"
$ws.Cells.Item(30, 13).Value = "***.*"
$ws.Cells.Item(30, 14).NumberFormat = "@"
$ws.Cells.Item(30, 14).Value = "***.*"
